$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 7).Value = 0.6848073333333332
$ws.Cells.Item(2, 8).Value = 2.054422
$ws.Cells.Item(2, 9).Value = 0.2268310526442471
$ws.Cells.Item(2, 10).Value = 0.2268310526442472
$ws.Cells.Item(2, 13).Value = 0.006825333333333333
$ws.Cells.Item(2, 14).Value = 0.020476
$ws.Cells.Item(2, 15).Value = 0.0506148883313352
$ws.Cells.Item(2, 16).Value = 0.0506148883313352
$ws.Cells.Item(2, 17).Value = 0.00467403831911111
$ws.Cells.Item(2, 18).Value = 0.042066344872
$ws.Cells.Item(2, 19).Value = 0.01148102839966778
$ws.Cells.Item(2, 20).Value = 0.01148102839966779

$ws.Cells.Item(3, 7).Value = 0.6848073333333332
$ws.Cells.Item(3, 8).Value = 2.054422
$ws.Cells.Item(3, 9).Value = 0.2268310526442471
$ws.Cells.Item(3, 10).Value = 0.2268310526442472
$ws.Cells.Item(3, 15).Value = 0.3671161428271267
$ws.Cells.Item(3, 16).Value = 0.3671161428271267
$ws.Cells.Item(3, 17).Value = 0.03390138703666666
$ws.Cells.Item(3, 18).Value = 0.30511248333
$ws.Cells.Item(3, 19).Value = 0.08327334112017293
$ws.Cells.Item(3, 20).Value = 0.08327334112017294

$ws.Cells.Item(4, 7).Value = 0.6848073333333332
$ws.Cells.Item(4, 8).Value = 2.054422
$ws.Cells.Item(4, 9).Value = 0.2268310526442471
$ws.Cells.Item(4, 10).Value = 0.2268310526442472
$ws.Cells.Item(4, 13).Value = 0.078518
$ws.Cells.Item(4, 14).Value = 0.235554
$ws.Cells.Item(4, 15).Value = 0.582268968841538
$ws.Cells.Item(4, 16).Value = 0.582268968841538
$ws.Cells.Item(4, 17).Value = 0.05376970219866666
$ws.Cells.Item(4, 18).Value = 0.483927319788
$ws.Cells.Item(4, 19).Value = 0.1320766831244064
$ws.Cells.Item(4, 20).Value = 0.1320766831244064

$ws.Cells.Item(5, 9).Value = 0.1086184939966157
$ws.Cells.Item(5, 10).Value = 0.1086184939966157
$ws.Cells.Item(5, 13).Value = 0.006825333333333333
$ws.Cells.Item(5, 14).Value = 0.020476
$ws.Cells.Item(5, 15).Value = 0.0506148883313352
$ws.Cells.Item(5, 16).Value = 0.0506148883313352
$ws.Cells.Item(5, 17).Value = 0.002238172407111111
$ws.Cells.Item(5, 18).Value = 0.020143551664
$ws.Cells.Item(5, 19).Value = 0.005497712944356505
$ws.Cells.Item(5, 20).Value = 0.005497712944356505

$ws.Cells.Item(6, 9).Value = 0.1086184939966157
$ws.Cells.Item(6, 10).Value = 0.1086184939966157
$ws.Cells.Item(6, 15).Value = 0.3671161428271267
$ws.Cells.Item(6, 16).Value = 0.3671161428271267
$ws.Cells.Item(6, 19).Value = 0.03987560255572897
$ws.Cells.Item(6, 20).Value = 0.03987560255572897

$ws.Cells.Item(7, 9).Value = 0.1086184939966157
$ws.Cells.Item(7, 10).Value = 0.1086184939966157
$ws.Cells.Item(7, 13).Value = 0.078518
$ws.Cells.Item(7, 14).Value = 0.235554
$ws.Cells.Item(7, 15).Value = 0.582268968841538
$ws.Cells.Item(7, 16).Value = 0.582268968841538
$ws.Cells.Item(7, 17).Value = 0.02574772725066667
$ws.Cells.Item(7, 18).Value = 0.231729545256
$ws.Cells.Item(7, 19).Value = 0.06324517849653019
$ws.Cells.Item(7, 20).Value = 0.06324517849653019

$ws.Cells.Item(8, 5).Value = 3
$ws.Cells.Item(8, 6).Value = 1
$ws.Cells.Item(8, 7).Value = 2.006290666666667
$ws.Cells.Item(8, 8).Value = 6.018872
$ws.Cells.Item(8, 9).Value = 0.6645504533591371
$ws.Cells.Item(8, 10).Value = 0.6645504533591372
$ws.Cells.Item(8, 13).Value = 0.006825333333333333
$ws.Cells.Item(8, 14).Value = 0.020476
$ws.Cells.Item(8, 15).Value = 0.0506148883313352
$ws.Cells.Item(8, 16).Value = 0.0506148883313352
$ws.Cells.Item(8, 17).Value = 0.01369360256355556
$ws.Cells.Item(8, 18).Value = 0.123242423072
$ws.Cells.Item(8, 19).Value = 0.03363614698731091
$ws.Cells.Item(8, 20).Value = 0.03363614698731091

$ws.Cells.Item(9, 5).Value = 3
$ws.Cells.Item(9, 6).Value = 1
$ws.Cells.Item(9, 7).Value = 2.006290666666667
$ws.Cells.Item(9, 8).Value = 6.018872
$ws.Cells.Item(9, 9).Value = 0.6645504533591371
$ws.Cells.Item(9, 10).Value = 0.6645504533591372
$ws.Cells.Item(9, 15).Value = 0.3671161428271267
$ws.Cells.Item(9, 16).Value = 0.3671161428271267
$ws.Cells.Item(9, 17).Value = 0.09932141945333334
$ws.Cells.Item(9, 18).Value = 0.8938927750800001
$ws.Cells.Item(9, 19).Value = 0.2439671991512248
$ws.Cells.Item(9, 20).Value = 0.2439671991512248

$ws.Cells.Item(10, 5).Value = 3
$ws.Cells.Item(10, 6).Value = 1
$ws.Cells.Item(10, 7).Value = 2.006290666666667
$ws.Cells.Item(10, 8).Value = 6.018872
$ws.Cells.Item(10, 9).Value = 0.6645504533591371
$ws.Cells.Item(10, 10).Value = 0.6645504533591372
$ws.Cells.Item(10, 13).Value = 0.078518
$ws.Cells.Item(10, 14).Value = 0.235554
$ws.Cells.Item(10, 15).Value = 0.582268968841538
$ws.Cells.Item(10, 16).Value = 0.582268968841538
$ws.Cells.Item(10, 17).Value = 0.1575299305653333
$ws.Cells.Item(10, 18).Value = 1.417769375088
$ws.Cells.Item(10, 19).Value = 0.3869471072206013
$ws.Cells.Item(10, 20).Value = 0.3869471072206014
